$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 284 - new section header "Project stage"
$ws.Range("A284").Value = "Project stage"

# Row 285 - Adjusted Benefits Cost Ratio (BCR) with three numeric values
$ws.Range("A285").Value = "Adjusted Benefits Cost Ratio (BCR)"
$ws.Range("B285").Value = 67
$ws.Range("C285").Value = 89
$ws.Range("D285").Value = 90

# Row 286 - Initial Benefits Cost Ratio (BCR)
$ws.Range("A286").Value = "Initial Benefits Cost Ratio (BCR)"

# Row 287 - VfM Category single entry
$ws.Range("A287").Value = "VfM Category single entry"

# Row 288 - VfM Category lower range
$ws.Range("A288").Value = "VfM Category lower range"

# Row 289 - VfM Category upper range
$ws.Range("A289").Value = "VfM Category upper range"

# Row 290 - SRO Benefits RAG with Green/Red/Green
$ws.Range("A290").Value = "SRO Benefits RAG"
$ws.Range("B290").Value = "Green"
$ws.Range("C290").Value = "Red"
$ws.Range("D290").Value = "Green"

# Update view: scroll/selection moved further down the sheet
$excel.ActiveWindow.ScrollRow = 272
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A291").Select()
